# "8 de junho - inesc"
# Adds a new reference row (R032 / Book / 2005 / "High-Speed Serial I/0 Made
# Simple" / SIM) to the bottom of the "Minhas Referencias" table on Sheet1,
# mirroring the layout of the existing rows (E:L merged for the title, M
# highlighted in yellow like the previous last row), then leaves the
# selection on D36 (just below the new row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 35 -----------------------------------------------------------
$ws.Range("A35").Value = "R032"
$ws.Range("B35").Value = "Book"
$ws.Range("C35").Value = "Book"
$ws.Range("D35").Value = 2005
$ws.Range("E35").Value = "High-Speed Serial I/0 Made Simple"

# Title spans E:L, same as every other row in the table.
$ws.Range("E35:L35").Merge() | Out-Null
$ws.Range("E35:L35").HorizontalAlignment = -4108
$ws.Range("E35:L35").VerticalAlignment = -4108

# "Corrigido" column - mark as done, highlighted like the row above it.
$ws.Range("M35").Value = "SIM"
$ws.Range("M35").Interior.Color = 65535
$ws.Range("M35").HorizontalAlignment = -4108
$ws.Range("M35").VerticalAlignment = -4108

# --- Final selection --------------------------------------------------
$ws.Range("D36").Select() | Out-Null
